{"js": "// Switch pandoc to nix.\n//\n// The underlying edit re-serializes a handful of run/style properties\n// (bold/italic toggles) so that `<w:b/>` / `<w:i/>` are written before\n// their `*Cs` / `color` siblings inside `<w:rPr>`. We reproduce that by\n// touching the Bold/Italic font properties on the exact runs (and\n// character styles) that changed \u2014 re-assigning the (already-true)\n// value is enough to make the host re-emit `rPr` in the new order.\n\nconst body = context.document.body;\n\n// --- document.xml: individual runs -----------------------------------\n// \"italics\" -> italic run, re-assert italic\nconst italicsResults = body.search(\"italics\", { matchCase: true });\nitalicsResults.load(\"items\");\nawait context.sync();\nitalicsResults.items[0].font.italic = true;\n\n// \"bold\" -> bold run, re-assert bold\nconst boldResults = body.search(\"bold\", { matchCase: true });\nboldResults.load(\"items\");\nawait context.sync();\nboldResults.items[0].font.bold = true;\n\n// \"line break\" appears twice; only the second (the formatted, standalone\n// run) carries italic styling -- match on it specifically.\nconst lineBreakResults = body.search(\"line break\", { matchCase: true });\nlineBreakResults.load(\"items\");\nawait context.sync();\nlineBreakResults.items[lineBreakResults.items.length - 1].font.italic = true;\n\n// \"formatting\" -> bold run, re-assert bold\nconst formattingResults = body.search(\"formatting\", { matchCase: true });\nformattingResults.load(\"items\");\nawait context.sync();\nformattingResults.items[0].font.bold = true;\n\n// \"even more italic text\" -> italic run, re-assert italic\nconst moreItalicResults = body.search(\"even more italic text\", { matchCase: true });\nmoreItalicResults.load(\"items\");\nawait context.sync();\nmoreItalicResults.items[0].font.italic = true;\n\nawait context.sync();\n\n// --- styles.xml: character styles -------------------------------------\nconst styles = context.document.getStyles();\n\nconst boldColorStyles = [\"KeywordTok\", \"ImportTok\", \"ControlFlowTok\", \"AlertTok\", \"ErrorTok\"];\nfor (const name of boldColorStyles) {\n  styles.getByName(name).font.bold = true;\n}\n\nconst italicColorStyles = [\"CommentTok\", \"DocumentationTok\"];\nfor (const name of italicColorStyles) {\n  styles.getByName(name).font.italic = true;\n}\n\nconst boldItalicColorStyles = [\"AnnotationTok\", \"CommentVarTok\", \"InformationTok\", \"WarningTok\"];\nfor (const name of boldItalicColorStyles) {\n  const s = styles.getByName(name);\n  s.font.bold = true;\n  s.font.italic = true;\n}\n\nawait context.sync();\n", "ps1": "# Switch pandoc to nix.\n#\n# The underlying edit re-serializes a handful of run/style properties\n# (bold/italic toggles) so that <w:b/> / <w:i/> are written before their\n# *Cs / color siblings inside <w:rPr>. We reproduce that by touching the\n# Bold/Italic font properties on the exact runs (and character styles)\n# that changed -- re-assigning the (already-true) value is enough to make\n# the host re-emit rPr in the new order.\n\n$d = $word.ActiveDocument\n\nfunction Find-NextRange($rng, $text) {\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $text\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 0\n    $rng.Find.Execute() | Out-Null\n}\n\n# --- document.xml: individual runs -----------------------------------\n\n# \"italics\" -> italic run, re-assert italic\n$rItalics = $d.Content\nFind-NextRange $rItalics \"italics\"\n$rItalics.Font.Italic = 1\n\n# \"bold\" -> bold run, re-assert bold\n$rBold = $d.Content\nFind-NextRange $rBold \"bold\"\n$rBold.Font.Bold = 1\n\n# \"line break\" appears twice; only the second (the formatted, standalone\n# run) carries italic styling -- advance past the first hit to reach it.\n$rLineBreak = $d.Content\nFind-NextRange $rLineBreak \"line break\"\n$rLineBreak.Collapse(0)\nFind-NextRange $rLineBreak \"line break\"\n$rLineBreak.Font.Italic = 1\n\n# \"formatting\" -> bold run, re-assert bold\n$rFormatting = $d.Content\nFind-NextRange $rFormatting \"formatting\"\n$rFormatting.Font.Bold = 1\n\n# \"even more italic text\" -> italic run, re-assert italic\n$rMoreItalic = $d.Content\nFind-NextRange $rMoreItalic \"even more italic text\"\n$rMoreItalic.Font.Italic = 1\n\n# --- styles.xml: character styles -------------------------------------\n\n$boldColorStyles = @(\"KeywordTok\", \"ImportTok\", \"ControlFlowTok\", \"AlertTok\", \"ErrorTok\")\nforeach ($name in $boldColorStyles) {\n    $s = $d.Styles($name)\n    $s.Font.Bold = 1\n}\n\n$italicColorStyles = @(\"CommentTok\", \"DocumentationTok\")\nforeach ($name in $italicColorStyles) {\n    $s = $d.Styles($name)\n    $s.Font.Italic = 1\n}\n\n$boldItalicColorStyles = @(\"AnnotationTok\", \"CommentVarTok\", \"InformationTok\", \"WarningTok\")\nforeach ($name in $boldItalicColorStyles) {\n    $s = $d.Styles($name)\n    $s.Font.Bold = 1\n    $s.Font.Italic = 1\n}\n"}
